# Debate, Waterfall 1, Waterfall 2
#
# Appends three new blank slides (same "Blank" layout as the deck's existing
# final slide) to the end of the presentation. Each new slide is produced by
# duplicating the previous slide in the chain, which keeps them in the
# correct final order (and mirrors how PowerPoint's own "Duplicate Slide"
# command behaves when used repeatedly).

$p = $ppt.ActivePresentation

$lastIndex = $p.Slides.Count
$sourceSlide = $p.Slides.Item($lastIndex)

# First new slide: duplicate of the existing last slide.
$range = $sourceSlide.Duplicate()
$newSlide = $range.Item(1)

# Second new slide: duplicate of the first new slide, so it lands right
# after it.
$range = $newSlide.Duplicate()
$newSlide = $range.Item(1)

# Third new slide: duplicate of the second new slide.
$range = $newSlide.Duplicate()
$newSlide = $range.Item(1)

Write-Host "Slide count is now $($p.Slides.Count)"
